# Jogos_da_Semana_FlashScore_2024-11-11.xlsx update
# - Row 8: refresh four odds values (O8:R8)
# - A brand-new match (Paraguay - Primera Division) is inserted as row 9,
#   pushing the former rows 9 and 10 down to rows 10 and 11
# - Row 11 (the match that used to live in row 10) also gets a handful of
#   its odds refreshed to newer prices

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update Row 8 odds ---
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 1.9
$ws.Range("R8").Value = 1.9

# --- Step 2: insert a new row at position 9 ---
# This shifts the old row 9 (Romania - Liga 1) down to row 10
# and the old row 10 (Spain - LaLiga2) down to row 11.
$ws.Rows.Item(9).Insert()

# --- Step 3: populate the new Row 9 with the Paraguay - Primera Division match ---
$ws.Range("A9").Value = "SlBBcd9a"

# Column B looks like a date ("11/11/2024"); force text so Excel doesn't
# convert it into a date serial number, then drop the temporary format so
# the cell ends up with no special style, just like its neighbours.
$dateCell = $ws.Range("B9")
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/11/2024"
$dateCell.ClearFormats()

$ws.Range("C9").Value = "20:00"
$ws.Range("D9").Value = "PARAGUAY - PRIMERA DIVISION"
$ws.Range("E9").Value = "Ameliano"
$ws.Range("F9").Value = "2 de Mayo"
$ws.Range("G9").Value = 2.38
$ws.Range("H9").Value = 2.8
$ws.Range("I9").Value = 3.25
$ws.Range("J9").Value = 3.4
$ws.Range("K9").Value = 1.83
$ws.Range("L9").Value = 4.33
$ws.Range("M9").Value = 1.14
$ws.Range("N9").Value = 5.5
$ws.Range("O9").Value = 1.57
$ws.Range("P9").Value = 2.25
$ws.Range("Q9").Value = 2.88
$ws.Range("R9").Value = 1.4
$ws.Range("S9").Value = 1.62
$ws.Range("T9").Value = 2.2
$ws.Range("U9").Value = 2.25
$ws.Range("V9").Value = 1.57
$ws.Range("W9").Value = 6
$ws.Range("X9").Value = 10
$ws.Range("Y9").Value = 11
$ws.Range("Z9").Value = 23
$ws.Range("AA9").Value = 26
$ws.Range("AB9").Value = 41
$ws.Range("AC9").Value = 5.5
$ws.Range("AD9").Value = 6
$ws.Range("AE9").Value = 21
$ws.Range("AF9").Value = 81
$ws.Range("AG9").Value = 1250
$ws.Range("AH9").Value = 7
$ws.Range("AI9").Value = 15
$ws.Range("AJ9").Value = 13
$ws.Range("AK9").Value = 34
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 51
$ws.Range("AN9").Value = 4.33
$ws.Range("AO9").Value = 15
$ws.Range("AP9").Value = 34
$ws.Range("AQ9").Value = 51
$ws.Range("AR9").Value = 101
$ws.Range("AS9").Value = 351
$ws.Range("AT9").Value = 2.2
$ws.Range("AU9").Value = 9.5
$ws.Range("AV9").Value = 81
$ws.Range("AW9").Value = 5
$ws.Range("AX9").Value = 21
$ws.Range("AY9").Value = 41
$ws.Range("AZ9").Value = 81
$ws.Range("BA9").Value = 126
$ws.Range("BB9").Value = 351
$ws.Range("BC9").Value = 51
$ws.Range("BD9").Value = 51

# --- Step 4: refresh odds on Row 11 (formerly Row 10 - Spain LaLiga2) ---
$ws.Range("G11").Value = 1.62
$ws.Range("H11").Value = 3.7
$ws.Range("I11").Value = 5.5
$ws.Range("J11").Value = 2.25
$ws.Range("Q11").Value = 1.95
$ws.Range("R11").Value = 1.9
$ws.Range("X11").Value = 7.5
$ws.Range("Z11").Value = 12
$ws.Range("AH11").Value = 15
$ws.Range("AI11").Value = 29
$ws.Range("AQ11").Value = 26
$ws.Range("AW11").Value = 7
$ws.Range("AX11").Value = 29
